# Translation sheet: remove the 10 "russian text" placeholder rows that were
# replaced by a single svg_1 asset, and re-center the alignment of the four
# rows that used to sit right before them.
#
# Net effect (matches the target OOXML diff):
#   - rows 23-30 (8 rows) are deleted, shifting rows 31+ up by 8
#   - the (now shifted) former rows 48-49 (2 more trailing rows) are deleted
#   - D15 / D17 / D19 / D21 ("Alignment" column) change from Left to Center

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Delete the block of 8 rows (originally rows 23..30).
$ws.Range("B23:F30").EntireRow.Delete()

# After the shift above, the rows that used to be 48 and 49 are now 40 and 41.
$ws.Range("B40:F41").EntireRow.Delete()

# Realign these four rows from Left to Center.
$ws.Range("D15").Value2 = "Center"
$ws.Range("D17").Value2 = "Center"
$ws.Range("D19").Value2 = "Center"
$ws.Range("D21").Value2 = "Center"
